$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: swap the B..AC contents (columns 2..29) of two rows, leaving
# column A (the running index) untouched.
# ---------------------------------------------------------------------------
function Swap-RowData {
    param($ws, $row1, $row2, $firstCol, $lastCol)

    $vals1 = @{}
    $vals2 = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals1[$c] = $ws.Cells.Item($row1, $c).Value2
        $vals2[$c] = $ws.Cells.Item($row2, $c).Value2
    }
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row1, $c).Value = $vals2[$c]
        $ws.Cells.Item($row2, $c).Value = $vals1[$c]
    }
}

# ---------------------------------------------------------------------------
# Helper: cyclically rotate the B..AC contents across an ordered list of
# rows, leaving column A untouched. Row at position i receives the old
# B..AC data of the row at position i+1 (wrapping around).
# newRows[i] = oldRows[i+1 (mod n)]
# ---------------------------------------------------------------------------
function Rotate-RowData {
    param($ws, $rows, $firstCol, $lastCol)

    $n = $rows.Count
    $snapshot = @()
    for ($i = 0; $i -lt $n; $i++) {
        $rowVals = @{}
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $rowVals[$c] = $ws.Cells.Item($rows[$i], $c).Value2
        }
        $snapshot += ,$rowVals
    }
    for ($i = 0; $i -lt $n; $i++) {
        $src = $snapshot[($i + 1) % $n]
        $destRow = $rows[$i]
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $ws.Cells.Item($destRow, $c).Value = $src[$c]
        }
    }
}

# Column B = 2 ... Column AC = 29
$firstCol = 2
$lastCol = 29

# Simple pairwise swaps
Swap-RowData $ws 64 65 $firstCol $lastCol
Swap-RowData $ws 71 72 $firstCol $lastCol
Swap-RowData $ws 95 96 $firstCol $lastCol
Swap-RowData $ws 120 121 $firstCol $lastCol

# 4-row cyclic rotation: new104<-old105, new105<-old107, new107<-old106, new106<-old104
# That is exactly the rotation where each row (in this listed order) takes
# the data of the "next" row in the list (wrapping around).
Rotate-RowData $ws @(104, 105, 107, 106) $firstCol $lastCol

# ---------------------------------------------------------------------------
# Append two new match rows (140 and 141) at the bottom of the sheet.
# ---------------------------------------------------------------------------

# Row 140 - full result already known
$ws.Range("A140").Value = 138
$ws.Range("A139").Copy()
$ws.Range("A140").PasteSpecial(-4122)

$ws.Range("B140").Value = 7721093
$ws.Range("C140").Value = "Estonia Meistriliiga"
$ws.Range("D140").Value = "Estonia Meistriliiga"
$ws.Range("E140").Value = 45399.54166666666
$ws.Range("E139").Copy()
$ws.Range("E140").PasteSpecial(-4122)
$ws.Range("F140").Value = "FC Flora Tallinn"
$ws.Range("G140").Value = "Paide Linnameeskond"
$ws.Range("H140").Value = 1
$ws.Range("I140").Value = 3
$ws.Range("J140").Value = "A"
$ws.Range("K140").Value = 2.3
$ws.Range("L140").Value = 3.4
$ws.Range("M140").Value = 2.6
$ws.Range("N140").Value = 2.4
$ws.Range("O140").Value = 3.5
$ws.Range("P140").Value = 2.4
$ws.Range("Q140").Value = 0
$ws.Range("R140").Value = 1.9
$ws.Range("S140").Value = 1.9
$ws.Range("T140").Value = 2.5
$ws.Range("U140").Value = 1.75
$ws.Range("V140").Value = 1.95
$ws.Range("W140").Value = -1
$ws.Range("X140").Value = -1
$ws.Range("Y140").Value = 1.4
$ws.Range("Z140").Value = -1
$ws.Range("AA140").Value = 0.8999999999999999
$ws.Range("AB140").Value = 0.75
$ws.Range("AC140").Value = -1

# Row 141 - match not played yet (no FTHG/FTAG/FTR, no PLH/PLD/PLA result columns)
$ws.Range("A141").Value = 139
$ws.Range("A139").Copy()
$ws.Range("A141").PasteSpecial(-4122)

$ws.Range("B141").Value = 7721089
$ws.Range("C141").Value = "Estonia Meistriliiga"
$ws.Range("D141").Value = "Estonia Meistriliiga"
$ws.Range("E141").Value = 45402.35416666666
$ws.Range("E141").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F141").Value = "Paide Linnameeskond"
$ws.Range("G141").Value = "JK Tammeka Tartu"
$ws.Range("K141").Value = 1.333
$ws.Range("L141").Value = 4.75
$ws.Range("M141").Value = 6.5
$ws.Range("N141").Value = 1.444
$ws.Range("O141").Value = 4.5
$ws.Range("P141").Value = 5
$ws.Range("Q141").Value = -1.25
$ws.Range("R141").Value = 1.95
$ws.Range("S141").Value = 1.85
$ws.Range("T141").Value = 3
$ws.Range("U141").Value = 2
$ws.Range("V141").Value = 1.8
$ws.Range("W141").Value = 0
$ws.Range("X141").Value = 0
$ws.Range("Y141").Value = 0
$ws.Range("Z141").Value = 0
$ws.Range("AA141").Value = 0
